$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

# 1. Fill in the date placeholder "June __, 2020" -> "June 30, 2020"
Replace-Text "__," "30,"

# 2. "...we present results from a surface sterilization treatment of invertebrate consumers prior to ..."
#    -> "...we present results of a study examining the effects of sterilization treatment of invertebrate consumers prior to ..."
Replace-Text "we present results from a surface sterilization treatment of invertebrate consumers prior to" "we present results of a study examining the effects of sterilization treatment of invertebrate consumers prior to"

# 3. " aimed at exploring the effects of surface contaminants on diet DNA data. " -> " on various metrics of diet. "
Replace-Text " aimed at exploring the effects of surface contaminants on diet DNA data. " " on various metrics of diet. "

# 4. "e believe this study directly relates to the scope of " -> "e believe this study is directly in the scope of "
Replace-Text "e believe this study directly relates to the scope of" "e believe this study is directly in the scope of"

# 5. "...most cases in validating diet DNA metabarcoding data from DNA extracted..."
#    -> "...most cases when using DNA metabarcoding data to gain diet information from DNA extracted..."
Replace-Text "sterilization does not appear to be a necessary step in most cases in validating diet DNA metabarcoding data from DNA extracted" "sterilization does not appear to be a necessary step in most cases when using DNA metabarcoding data to gain diet information from DNA extracted"

# 6. "...DNA metabarcoding continues to provide valuable insight into a range of consumptive interactions..."
#    -> "...DNA metabarcoding is increasingly used to identify consumptive interactions..."
Replace-Text "DNA metabarcoding continues to provide valuable insight into a range of consumptive interactions" "DNA metabarcoding is increasingly used to identify consumptive interactions"

# 7. " environmental contamination and systematic fixes (e.g. surface sterilization). Our study provides validation
#    for past and future studies that use this approach and provides suggestions..."
#    -> " environmental contamination and there is no consensus on best practice for use of surface sterilization
#    to address these risks. Our study provides validation for past and future studies that do not use surface
#    sterilization and provides suggestions..."
Replace-Text " environmental contamination and systematic fixes (e.g. surface sterilization). Our study provides validation for past and future studies that use this approach and provides suggestions" " environmental contamination and there is no consensus on best practice for use of surface sterilization to address these risks. Our study provides validation for past and future studies that do not use surface sterilization and provides suggestions"

# 8. "will continue to expand the utility of this method in  large-scale ecological questions, such as the
#    maintenance of biodiversity and ecosystem functions.   "
#    -> "will be critical to improving conclusions from these types of studies and facilitating cross-study
#    comparison.   "
Replace-Text "will continue to expand the utility of this method in  large-scale ecological questions, such as the maintenance of biodiversity and ecosystem functions.   " "will be critical to improving conclusions from these types of studies and facilitating cross-study comparison.  "
